$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22284.234
$ws.Range("I62").Value = 33789.6
$ws.Range("K62").Value = 33789.6
$ws.Range("M62").Value = -33165.6
$ws.Range("H65").Value = 22284.234
$ws.Range("I65").Value = 33789.6
$ws.Range("K65").Value = 168948
$ws.Range("M65").Value = -165828
$ws.Range("H70").Value = 2974.4666
$ws.Range("I70").Value = 2141.3333
$ws.Range("J70").Value = 3529.889
$ws.Range("K70").Value = 6423.999899999999
$ws.Range("L70").Value = 10589.667
$ws.Range("M70").Value = -6153.999899999999
$ws.Range("N70").Value = -11129.667
$ws.Range("H73").Value = 2974.4666
$ws.Range("I73").Value = 2141.3333
$ws.Range("J73").Value = 3529.889
$ws.Range("K73").Value = 6423.999899999999
$ws.Range("L73").Value = 10589.667
$ws.Range("M73").Value = -5487.999899999999
$ws.Range("N73").Value = -12461.667
$ws.Range("H106").Value = 7653.0625
$ws.Range("I106").Value = 5831.25
$ws.Range("K106").Value = 5831.25
$ws.Range("M106").Value = -5200.25
$ws.Range("H107").Value = 2152.0908
$ws.Range("J107").Value = 1965.8334
$ws.Range("L107").Value = 1965.8334
$ws.Range("N107").Value = -5805.8334
$ws.Range("H116").Value = 5084.8
$ws.Range("J116").Value = 5674
$ws.Range("L116").Value = 5674
$ws.Range("N116").Value = -12558
$ws.Range("H132").Value = 4363.755
$ws.Range("I132").Value = 3003.1365
$ws.Range("K132").Value = 9009.4095
$ws.Range("M132").Value = -6479.4095
$ws.Range("H138").Value = 2850.1628
$ws.Range("I138").Value = 1253.6923
$ws.Range("J138").Value = 5291.8237
$ws.Range("K138").Value = 3761.0769
$ws.Range("L138").Value = 15875.4711
$ws.Range("M138").Value = 1378.9231
$ws.Range("N138").Value = -26155.4711
$ws.Range("H141").Value = 4883.125
$ws.Range("I141").Value = 4883.125
$ws.Range("K141").Value = 14649.375
$ws.Range("M141").Value = -9469.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2063.182
$ws.Range("J2").Value = 3128
$ws.Range("L2").Value = 3128
$ws.Range("N2").Value = -3354
$ws.Range("H32").Value = 19147.047
$ws.Range("I32").Value = 18909.855
$ws.Range("K32").Value = 18909.855
$ws.Range("M32").Value = -18622.855
$ws.Range("H43").Value = 19998
$ws.Range("I43").Value = 19998
$ws.Range("K43").Value = 19998
$ws.Range("M43").Value = -19685
$ws.Range("H63").Value = 4900.6665
$ws.Range("I63").Value = 1801.3334
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 1801.3334
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -1115.3334
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 4900.6665
$ws.Range("I66").Value = 1801.3334
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 9006.666999999999
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -5574.666999999999
$ws.Range("N66").Value = -46864
$ws.Range("H116").Value = 2063.182
$ws.Range("J116").Value = 3128
$ws.Range("L116").Value = 3128
$ws.Range("N116").Value = -7716
$ws.Range("H122").Value = 3805.375
$ws.Range("I122").Value = 3675.5
$ws.Range("J122").Value = 4195
$ws.Range("K122").Value = 11026.5
$ws.Range("L122").Value = 12585
$ws.Range("M122").Value = -8576.5
$ws.Range("N122").Value = -17485
$ws.Range("H131").Value = 92882.336
$ws.Range("J131").Value = 97999.5
$ws.Range("L131").Value = 97999.5
$ws.Range("N131").Value = -108079.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2063.182
$ws.Range("J3").Value = 3128
$ws.Range("L3").Value = 3128
$ws.Range("N3").Value = -3356

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2585.739
$ws.Range("I31").Value = 2041.6111
$ws.Range("J31").Value = 4544.6
$ws.Range("K31").Value = 2041.6111
$ws.Range("L31").Value = 4544.6
$ws.Range("M31").Value = -1746.6111
$ws.Range("N31").Value = -5134.6
$ws.Range("H34").Value = 2585.739
$ws.Range("I34").Value = 2041.6111
$ws.Range("J34").Value = 4544.6
$ws.Range("K34").Value = 2041.6111
$ws.Range("L34").Value = 4544.6
$ws.Range("M34").Value = -1839.6111
$ws.Range("N34").Value = -4948.6
$ws.Range("H99").Value = 9147.929
$ws.Range("I99").Value = 11544.3
$ws.Range("K99").Value = 11544.3
$ws.Range("M99").Value = -10046.3
$ws.Range("H105").Value = 50749
$ws.Range("I105").Value = 66998.664
$ws.Range("K105").Value = 66998.664
$ws.Range("M105").Value = -65251.664
$ws.Range("H126").Value = 9147.929
$ws.Range("I126").Value = 11544.3
$ws.Range("K126").Value = 34632.89999999999
$ws.Range("M126").Value = -32162.89999999999
$ws.Range("H132").Value = 3314.6
$ws.Range("I132").Value = 3225.2307
$ws.Range("K132").Value = 9675.6921
$ws.Range("M132").Value = -7145.6921
$ws.Range("H134").Value = 32439.121
$ws.Range("I134").Value = 38648.074
$ws.Range("J134").Value = 4498.8335
$ws.Range("K134").Value = 115944.222
$ws.Range("L134").Value = 13496.5005
$ws.Range("M134").Value = -113409.222
$ws.Range("N134").Value = -18566.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1874415.9
$ws.Range("I4").Value = 2122395
$ws.Range("K4").Value = 6367185
$ws.Range("M4").Value = -6367073
$ws.Range("H59").Value = 3849.8333
$ws.Range("I59").Value = 2420
$ws.Range("J59").Value = 10999
$ws.Range("K59").Value = 7260
$ws.Range("L59").Value = 32997
$ws.Range("M59").Value = -6720
$ws.Range("N59").Value = -34077

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H70").Value = 11060.728
$ws.Range("H73").Value = 11060.728
$ws.Range("H92").Value = 37633.332
$ws.Range("J92").Value = 37633.332
$ws.Range("L92").Value = 37633.332
$ws.Range("N92").Value = -41377.332
$ws.Range("H96").Value = 9260.5
$ws.Range("I96").Value = 10260
$ws.Range("J96").Value = 8261
$ws.Range("K96").Value = 10260
$ws.Range("L96").Value = 8261
$ws.Range("M96").Value = -7514
$ws.Range("N96").Value = -13753
$ws.Range("H102").Value = 6322
$ws.Range("I102").Value = 1071.7333
$ws.Range("J102").Value = 26010.5
$ws.Range("K102").Value = 1071.7333
$ws.Range("L102").Value = 26010.5
$ws.Range("M102").Value = 550.2666999999999
$ws.Range("N102").Value = -29254.5
$ws.Range("H113").Value = 107571.65
$ws.Range("I113").Value = 71214.56
$ws.Range("K113").Value = 71214.56
$ws.Range("M113").Value = -69044.56
$ws.Range("H126").Value = 6137.8945
$ws.Range("I126").Value = 6240.385
$ws.Range("K126").Value = 18721.155
$ws.Range("M126").Value = -16251.155
$ws.Range("H132").Value = 45935.39
$ws.Range("I132").Value = 52304.35
$ws.Range("K132").Value = 156913.05
$ws.Range("M132").Value = -154383.05

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5460.276
$ws.Range("I40").Value = 5937.8076
$ws.Range("J40").Value = 1321.6666
$ws.Range("K40").Value = 5937.8076
$ws.Range("L40").Value = 1321.6666
$ws.Range("M40").Value = -5801.8076
$ws.Range("N40").Value = -1593.6666
$ws.Range("H122").Value = 424206.66
$ws.Range("I122").Value = 11926.9
$ws.Range("J122").Value = 718692.2
$ws.Range("K122").Value = 35780.7
$ws.Range("L122").Value = 2156076.6
$ws.Range("M122").Value = -33330.7
$ws.Range("N122").Value = -2160976.6
$ws.Range("H136").Value = 3320.318
$ws.Range("I136").Value = 2943.7646
$ws.Range("K136").Value = 8831.293799999999
$ws.Range("M136").Value = -6281.293799999999
$ws.Range("H137").Value = 88189.5
$ws.Range("J137").Value = 88189.5
$ws.Range("L137").Value = 88189.5
$ws.Range("N137").Value = -98389.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 25546.666
$ws.Range("I4").Value = 25546.666
$ws.Range("K4").Value = 25546.666
$ws.Range("M4").Value = -25433.666
$ws.Range("H28").Value = 30000
$ws.Range("J28").Value = 30000
$ws.Range("L28").Value = 30000
$ws.Range("N28").Value = -30696
$ws.Range("H47").Value = 3250
$ws.Range("I47").Value = 3000
$ws.Range("J47").Value = 4166.6665
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 4166.6665
$ws.Range("M47").Value = -2428
$ws.Range("N47").Value = -5310.6665
$ws.Range("H58").Value = 29308
$ws.Range("I58").Value = 29846.5
$ws.Range("J58").Value = 25000
$ws.Range("K58").Value = 29846.5
$ws.Range("L58").Value = 25000
$ws.Range("M58").Value = -29538.5
$ws.Range("N58").Value = -25616
$ws.Range("H96").Value = 3659.8
$ws.Range("I96").Value = 3150
$ws.Range("K96").Value = 3150
$ws.Range("M96").Value = -1777
$ws.Range("H100").Value = 374.14285
$ws.Range("I100").Value = 303.66666
$ws.Range("K100").Value = 607.33332
$ws.Range("M100").Value = -66.33331999999996
$ws.Range("H122").Value = 113701.22
$ws.Range("I122").Value = 2384.5
$ws.Range("K122").Value = 7153.5
$ws.Range("M122").Value = -4703.5
$ws.Range("H130").Value = 44999.668
$ws.Range("J130").Value = 44999.668
$ws.Range("L130").Value = 44999.668
$ws.Range("N130").Value = -55039.668
$ws.Range("H132").Value = 23572.12
$ws.Range("I132").Value = 25591.334
$ws.Range("K132").Value = 76774.00199999999
$ws.Range("M132").Value = -74244.00199999999

Write-Output "Applied all changes"